# Generate Report for Handoff
#
# Regenerates the localization-status report: a new source file
# (16f4bed0-8b99-4d3d-ac91-84247a6c6002.md) replaces the old one
# (77ef3fa6-132a-4d52-8739-be415dae2553.md) everywhere it is referenced,
# the two per-locale handoff .xlf file names pick up their new content
# hashes, and the "latest" timestamps advance a few seconds.

$wb = $excel.ActiveWorkbook

$oldBase = "77ef3fa6-132a-4d52-8739-be415dae2553"
$newBase = "16f4bed0-8b99-4d3d-ac91-84247a6c6002"

$oldMd = "$oldBase.md"
$newMd = "$newBase.md"

$oldZhXlf = "$oldBase.b5c6d2550c570577c5c6004d90f2f8845cb357d0.zh-cn.xlf"
$newZhXlf = "$newBase.b95ab77519e70a783ed97edf66d60250b2027203.zh-cn.xlf"

$oldDeXlf = "$oldBase.b5c6d2550c570577c5c6004d90f2f8845cb357d0.de-de.xlf"
$newDeXlf = "$newBase.b95ab77519e70a783ed97edf66d60250b2027203.de-de.xlf"

$oldHoDate = "2016-09-05 23:10:51"
$newHoDate = "2016-09-05 23:11:15"

$oldZhDate = "2016-09-05 23:10:46"
$newZhDate = "2016-09-05 23:11:11"

# All three sheets' A2/B2 hyperlinks point at the same GitHub blob URL.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79560849d0d7d9ccb07ba126d5fc9d78546975db/e2e/77ef3fa6-132a-4d52-8739-be415dae2553.md"
$newOverviewDisplay = "e2e\" + $newMd

# Helper: update the display text of a single hyperlinked cell while
# keeping its existing target/relationship (Address) untouched. (Reading
# Hyperlink.Address back from the engine yields "", so the known external
# target is passed in explicitly instead of round-tripping it.)
function Update-HyperlinkCell($range, $address, $newText) {
    $range.Value = $newText
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newText) | Out-Null
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
Update-HyperlinkCell $wsOverview.Range("B2") $hyperlinkAddress $newOverviewDisplay
$wsOverview.Range("G2").Value = $newHoDate

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HyperlinkCell $wsZhCn.Range("A2") $hyperlinkAddress $newMd
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HyperlinkCell $wsDeDe.Range("A2") $hyperlinkAddress $newMd
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHoDate
